$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 87, pushing the existing rows 87..156 down
# to 88..157 (dimension grows from A1:T156 to A1:T157).
$ws.Rows.Item(87).Insert()

# Populate the newly inserted row 87 with the new record.
$ws.Range("A87").Value = 1
$ws.Range("B87").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C87").Value = "Arica y Parinacota"
$ws.Range("D87").Value = 45126
$ws.Range("E87").Value = 15
$ws.Range("F87").Value = "Fruta"
$ws.Range("G87").Value = 100102
$ws.Range("H87").Value = "Cítricos"
$ws.Range("I87").Value = 100102005
$ws.Range("J87").Value = "Naranja"
$ws.Range("K87").Value = "Fukumoto"
$ws.Range("L87").Value = "Segunda"
$ws.Range("M87").Value = 285
$ws.Range("N87").Value = 750
$ws.Range("O87").Value = 900
$ws.Range("P87").Value = 821
$ws.Range("Q87").Value = "$/kilo (en caja de 20 kilos)"
$ws.Range("R87").Value = "Región de Coquimbo"
$ws.Range("S87").Value = 821
$ws.Range("T87").Value = 1
